$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44692
$ws.Range("J2").Value = 120

$ws.Range("D3").Value = 44687
$ws.Range("J3").Value = 160

$ws.Range("D5").Value = 44691
$ws.Range("J5").Value = 100
